# Add variables for critical care beds (icu/ccu/bicu/sicu/othspec beds)
# to the "Lookup Table" and "Type and Label" sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Lookup Table")
$ws2 = $wb.Worksheets.Item("Type and Label")

# ---------------------------------------------------------------------
# Block A ("Lookup Table"): insert 5 rows after row 28 (fmt=10 section),
# before the existing "beds_total" (fmt=10) row.
# ---------------------------------------------------------------------
$ws1.Rows("29:33").Insert()

# Filled column-by-column (A29:A33, then B29:B33, ...) to mirror how the
# original edit was typed (also matches the shared-string insertion order
# recorded in the target workbook).
$colA = @("icu_beds", "ccu_beds", "bicu_beds", "sicu_beds", "othspec_beds")
$colB = @("S300001", "S300001", "S300001", "S300001", "S300001")
$colC = @("00200", "00200", "00200", "00200", "00200")
$colD = @("00800", "00900", "01000", "01100", "01200")
$colE = @(10, 10, 10, 10, 10)
$colF = @(1, 1, 1, 1, 1)

for ($i = 0; $i -lt 5; $i++) { $ws1.Range("A" + (29 + $i)).Value = $colA[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("B" + (29 + $i)).Value = $colB[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("C" + (29 + $i)).Value = $colC[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("D" + (29 + $i)).Value = $colD[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("E" + (29 + $i)).Value = $colE[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("F" + (29 + $i)).Value = $colF[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("G" + (29 + $i)).NumberFormat = "@" }

# ---------------------------------------------------------------------
# "Type and Label": insert 5 rows before the existing "beds_total" row
# (old row 29) and give each new rec its type + label.
# ---------------------------------------------------------------------
$ws2.Rows("29:33").Insert()

$colA2 = @("icu_beds", "ccu_beds", "bicu_beds", "sicu_beds", "othspec_beds")
$colB2 = @("stock", "stock", "stock", "stock", "stock")
$colC2 = @("intensive care unit beds", "coronary care unit beds", `
    "burn intensive care unit beds", "surgical intensive care unit beds", `
    "other special care beds")

for ($i = 0; $i -lt 5; $i++) { $ws2.Range("A" + (29 + $i)).Value = $colA2[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws2.Range("B" + (29 + $i)).Value = $colB2[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws2.Range("C" + (29 + $i)).Value = $colC2[$i] }

$ws2.Activate()
$ws2.Range("C39").Select()

# ---------------------------------------------------------------------
# Block B ("Lookup Table"): insert 5 rows (fmt=96 section), before the
# existing "beds_total" (fmt=96) row. After block A's insert, that row
# is now at row 64 (was row 59, +5 shift).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Rows("64:68").Insert()

$colA3 = @("icu_beds", "ccu_beds", "bicu_beds", "sicu_beds", "othspec_beds")
$colB3 = @("S300001", "S300001", "S300001", "S300001", "S300001")
$colC3 = @("0100", "0100", "0100", "0100", "0100")
$colD3 = @("02600", "02700", "02800", "02900", "02140")
$colE3 = @(96, 96, 96, 96, 96)
$colF3 = @(1, 1, 1, 1, 1)

for ($i = 0; $i -lt 5; $i++) { $ws1.Range("A" + (64 + $i)).Value = $colA3[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("B" + (64 + $i)).Value = $colB3[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("C" + (64 + $i)).Value = $colC3[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("D" + (64 + $i)).Value = $colD3[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("E" + (64 + $i)).Value = $colE3[$i] }
for ($i = 0; $i -lt 5; $i++) { $ws1.Range("F" + (64 + $i)).Value = $colF3[$i] }
# Rows 65-68 (not 64) carry the leftover empty styled G cell, matching
# the source edit.
$ws1.Range("G65").NumberFormat = "@"
$ws1.Range("G66").NumberFormat = "@"
$ws1.Range("G67").NumberFormat = "@"
$ws1.Range("G68").NumberFormat = "@"

$ws1.Activate()
$ws1.Range("F63:F68").Select()
